$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.139517664909363
$ws.Range("B1").Value = 2.596094369888306
$ws.Range("C1").Value = 6.093301773071289
$ws.Range("D1").Value = 2.143232822418213
$ws.Range("E1").Value = 1.234320878982544
